$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A20").Font.Name = "BCSans-Regular"
$ws.Range("A20").Font.Size = 18
$ws.Range("A20").Font.ThemeColor = 2
Write-Host "ok"
